# Scheduled market-data refresh for the Ramuh_Profits leve-crafting workbook.
# Each worksheet (ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR) is a per-job "Table_<job>"
# of turn-in leves; columns H-N hold the scraped/derived market figures:
#   H currentAveragePrice    I currentAveragePriceNQ   J currentAveragePriceHQ
#   K LevePriceNQ            L LevePriceHQ
#   M LeveProfitNQ           N LeveProfitHQ
# The runner refreshes these from the latest market-board pull; everything
# else on each row (leve name, item, level, exp, gil, ids) is untouched.

$wb = $excel.ActiveWorkbook

# ---- ALC -------------------------------------------------------------
$ws = $wb.Worksheets.Item("ALC")

# Row 12: Don't Be So Tallow / Beeswax
$ws.Range("H12").Value = 50000284
$ws.Range("I12").Value = 178.75
$ws.Range("J12").Value = 83333690
$ws.Range("K12").Value = 178.75
$ws.Range("L12").Value = 83333690
$ws.Range("M12").Value = -8.75
$ws.Range("N12").Value = -83334030

# Row 40: Stuck in the Moment / Horn Glue
$ws.Range("H40").Value = 2131.2
$ws.Range("I40").Value = 1226
$ws.Range("J40").Value = 3036.4
$ws.Range("K40").Value = 1226
$ws.Range("L40").Value = 3036.4
$ws.Range("M40").Value = -1051
$ws.Range("N40").Value = -3386.4

# Row 64: Forged from the Void / Void Glue
$ws.Range("H64").Value = 62503876
$ws.Range("I64").Value = 125002070
$ws.Range("J64").Value = 5675
$ws.Range("K64").Value = 125002070
$ws.Range("L64").Value = 5675
$ws.Range("M64").Value = -125001822
$ws.Range("N64").Value = -6171

# Row 67: Dodging the Draft (L) / Void Glue
$ws.Range("H67").Value = 62503876
$ws.Range("I67").Value = 125002070
$ws.Range("J67").Value = 5675
$ws.Range("K67").Value = 125002070
$ws.Range("L67").Value = 5675
$ws.Range("M67").Value = -125001212
$ws.Range("N67").Value = -7391

# Row 76: Warding Off Temptation / Enchanted Hardsilver Ink
$ws.Range("H76").Value = 5060.5713
$ws.Range("I76").Value = 2878.6
$ws.Range("J76").Value = 6272.778
$ws.Range("K76").Value = 2878.6
$ws.Range("L76").Value = 6272.778
$ws.Range("M76").Value = -2563.6
$ws.Range("N76").Value = -6902.778

# Row 79: The Garden of Arcane Delights (L) / Enchanted Hardsilver Ink
$ws.Range("H79").Value = 5060.5713
$ws.Range("I79").Value = 2878.6
$ws.Range("J79").Value = 6272.778
$ws.Range("K79").Value = 2878.6
$ws.Range("L79").Value = 6272.778
$ws.Range("M79").Value = -1786.6
$ws.Range("N79").Value = -8456.778

# Row 120: Supreme Official Strategy Guide / Dwarven Mythril Codex
# (HQ price now resolvable -> LeveProfitHQ column appears for the first time)
$ws.Range("H120").Value = 48000
$ws.Range("J120").Value = 48000
$ws.Range("L120").Value = 48000
$ws.Range("N120").Value = -57676

# Row 124: Luncheon Bound / Luncheon Toadskin Codex
$ws.Range("H124").Value = 41862.5
$ws.Range("J124").Value = 41862.5
$ws.Range("L124").Value = 41862.5
$ws.Range("N124").Value = -51682.5

# ---- ARM -------------------------------------------------------------
$ws = $wb.Worksheets.Item("ARM")

# Row 2: Ain't Got No Ingots / Bronze Ingot
$ws.Range("H2").Value = 772.913
$ws.Range("I2").Value = 621.17145
$ws.Range("J2").Value = 1255.7273
$ws.Range("K2").Value = 621.17145
$ws.Range("L2").Value = 1255.7273
$ws.Range("M2").Value = -508.17145
$ws.Range("N2").Value = -1481.7273

# Row 63: Rivets Run through It / Mythrite Rivets
$ws.Range("H63").Value = 2225
$ws.Range("I63").Value = 2100
$ws.Range("J63").Value = 2600
$ws.Range("K63").Value = 2100
$ws.Range("L63").Value = 2600
$ws.Range("M63").Value = -1414
$ws.Range("N63").Value = -3972

# Row 66: A Riveting Revival (L) / Mythrite Rivets
$ws.Range("H66").Value = 2225
$ws.Range("I66").Value = 2100
$ws.Range("J66").Value = 2600
$ws.Range("K66").Value = 10500
$ws.Range("L66").Value = 13000
$ws.Range("M66").Value = -7068
$ws.Range("N66").Value = -19864

# Row 88: The Mast Chance / Adamantite Rivets
$ws.Range("H88").Value = 2717.4348
$ws.Range("I88").Value = 2277.9443
$ws.Range("J88").Value = 4299.6
$ws.Range("K88").Value = 2277.9443
$ws.Range("L88").Value = 4299.6
$ws.Range("M88").Value = -1871.9443
$ws.Range("N88").Value = -5111.6

# Row 91: The Rose and the Riveter (L) / Adamantite Rivets
$ws.Range("H91").Value = 2717.4348
$ws.Range("I91").Value = 2277.9443
$ws.Range("J91").Value = 4299.6
$ws.Range("K91").Value = 2277.9443
$ws.Range("L91").Value = 4299.6
$ws.Range("M91").Value = -873.9443000000001
$ws.Range("N91").Value = -7107.6

# Row 102: Smells of Rich Tama-hagane / Tama-hagane Ingot
$ws.Range("H102").Value = 1982.2
$ws.Range("I102").Value = 1800
$ws.Range("J102").Value = 2027.75
$ws.Range("K102").Value = 1800
$ws.Range("L102").Value = 2027.75
$ws.Range("M102").Value = -178
$ws.Range("N102").Value = -5271.75

# Row 116: No Scope / Titanbronze Ingot
$ws.Range("H116").Value = 772.913
$ws.Range("I116").Value = 621.17145
$ws.Range("J116").Value = 1255.7273
$ws.Range("K116").Value = 621.17145
$ws.Range("L116").Value = 1255.7273
$ws.Range("M116").Value = 1672.82855
$ws.Range("N116").Value = -5843.7273

# Row 122: Haste for High Durium / High Durium Nugget
$ws.Range("H122").Value = 2037.7273
$ws.Range("I122").Value = 2263
$ws.Range("J122").Value = 1953.25
$ws.Range("K122").Value = 6789
$ws.Range("L122").Value = 5859.75
$ws.Range("M122").Value = -4339
$ws.Range("N122").Value = -10759.75

# ---- BSM -------------------------------------------------------------
$ws = $wb.Worksheets.Item("BSM")

# Row 3: Hells Bells / Bronze Ingot
$ws.Range("H3").Value = 772.913
$ws.Range("I3").Value = 621.17145
$ws.Range("J3").Value = 1255.7273
$ws.Range("K3").Value = 621.17145
$ws.Range("L3").Value = 1255.7273
$ws.Range("M3").Value = -507.17145
$ws.Range("N3").Value = -1483.7273

# Row 103: The Bigger the Blade / Doman Steel Tachi
# (HQ price no longer resolvable -> LeveProfitHQ column drops off entirely)
$ws.Range("H103").Value = 0
$ws.Range("J103").Value = 0
$ws.Range("L103").Value = 0
$ws.Range("N103").ClearContents()

# Row 105: Ingot to Wing It / Molybdenum Ingot
$ws.Range("H105").Value = 1832.3684
$ws.Range("I105").Value = 1336.125
$ws.Range("J105").Value = 2193.2727
$ws.Range("K105").Value = 1336.125
$ws.Range("L105").Value = 2193.2727
$ws.Range("M105").Value = 410.875
$ws.Range("N105").Value = -5687.2727

# Row 107: The Gold Experience / Deepgold Nugget
$ws.Range("H107").Value = 2070.484
$ws.Range("I107").Value = 1773
$ws.Range("J107").Value = 2611.3635
$ws.Range("K107").Value = 1773
$ws.Range("L107").Value = 2611.3635
$ws.Range("M107").Value = 147
$ws.Range("N107").Value = -6451.363499999999

# ---- CRP -------------------------------------------------------------
$ws = $wb.Worksheets.Item("CRP")

# Row 62: Splinter in the Sewers / Cedar Lumber
$ws.Range("H62").Value = 4870.8335
$ws.Range("I62").Value = 5527.778
$ws.Range("J62").Value = 2900
$ws.Range("K62").Value = 5527.778
$ws.Range("L62").Value = 2900
$ws.Range("M62").Value = -4903.778
$ws.Range("N62").Value = -4148

# Row 65: The Lumber of Their Discontent (L) / Cedar Lumber
$ws.Range("H65").Value = 4870.8335
$ws.Range("I65").Value = 5527.778
$ws.Range("J65").Value = 2900
$ws.Range("K65").Value = 27638.89
$ws.Range("L65").Value = 14500
$ws.Range("M65").Value = -24518.89
$ws.Range("N65").Value = -20740

# Row 111: Taking Aim / Applewood Longbow
# (HQ price no longer resolvable -> LeveProfitHQ column drops off entirely)
$ws.Range("H111").Value = 0
$ws.Range("J111").Value = 0
$ws.Range("L111").Value = 0
$ws.Range("N111").ClearContents()

# ---- CUL -------------------------------------------------------------
$ws = $wb.Worksheets.Item("CUL")

# Row 131: The Mountain Steeped / Tsai tou Vounou
$ws.Range("H131").Value = 861.0599999999999
$ws.Range("J131").Value = 956.6047
$ws.Range("L131").Value = 2869.8141
$ws.Range("N131").Value = -12949.8141

# ---- GSM -------------------------------------------------------------
$ws = $wb.Worksheets.Item("GSM")

# Row 14: All That Glitters / Copper Ear Cuffs
$ws.Range("H14").Value = 175243.44
$ws.Range("I14").Value = 223402.72
$ws.Range("J14").Value = 1870
$ws.Range("K14").Value = 223402.72
$ws.Range("L14").Value = 1870
$ws.Range("M14").Value = -223234.72
$ws.Range("N14").Value = -2206

# Row 70: Sky Is the Limit / Mythrite Ingot
$ws.Range("H70").Value = 4001
$ws.Range("I70").Value = 4001.4
$ws.Range("J70").Value = 4000
$ws.Range("K70").Value = 4001.4
$ws.Range("L70").Value = 4000
$ws.Range("M70").Value = -3731.4
$ws.Range("N70").Value = -4540

# Row 73: Hulls of Broken Dreams (L) / Mythrite Ingot
$ws.Range("H73").Value = 4001
$ws.Range("I73").Value = 4001.4
$ws.Range("J73").Value = 4000
$ws.Range("K73").Value = 4001.4
$ws.Range("L73").Value = 4000
$ws.Range("M73").Value = -3065.4
$ws.Range("N73").Value = -5872

# Row 80: Needs More Prayerbell / Hardsilver Ingot
$ws.Range("H80").Value = 3093.5
$ws.Range("I80").Value = 2227.5
$ws.Range("K80").Value = 2227.5
$ws.Range("M80").Value = -1229.5

# Row 83: With a Noise That Reaches Heaven (L) / Hardsilver Ingot
$ws.Range("H83").Value = 3093.5
$ws.Range("I83").Value = 2227.5
$ws.Range("K83").Value = 11137.5
$ws.Range("M83").Value = -6145.5

# Row 102: Put the Metal to the Peddle / Durium Ingot
$ws.Range("H102").Value = 3402.3704
$ws.Range("I102").Value = 3860.4443
$ws.Range("J102").Value = 2486.2222
$ws.Range("K102").Value = 3860.4443
$ws.Range("L102").Value = 2486.2222
$ws.Range("M102").Value = -2238.4443
$ws.Range("N102").Value = -5730.2222

# Row 113: Copious Crystal Cannons / Manasilver Nugget
$ws.Range("H113").Value = 626041.1
$ws.Range("I113").Value = 1111855.5
$ws.Range("K113").Value = 1111855.5
$ws.Range("M113").Value = -1109685.5

# Row 122: Awarding Academic Excellence / Ametrine
$ws.Range("H122").Value = 2818.4546
$ws.Range("I122").Value = 1160.8572
$ws.Range("K122").Value = 3482.5716
$ws.Range("M122").Value = -1032.5716

# ---- LTW -------------------------------------------------------------
$ws = $wb.Worksheets.Item("LTW")

# Row 40: Best Served Toad / Toad Leather
$ws.Range("H40").Value = 2267.3333
$ws.Range("I40").Value = 2147.762
$ws.Range("J40").Value = 2685.8333
$ws.Range("K40").Value = 2147.762
$ws.Range("L40").Value = 2685.8333
$ws.Range("M40").Value = -2011.762
$ws.Range("N40").Value = -2957.8333

# ---- WVR -------------------------------------------------------------
$ws = $wb.Worksheets.Item("WVR")

# Row 123: Helping Handwear / Fingerless Darkhempen Gloves of Healing
# (NQ price no longer resolvable -> LeveProfitNQ column drops off entirely)
$ws.Range("H123").Value = 35715.5
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 35715.5
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 35715.5
$ws.Range("M123").ClearContents()
$ws.Range("N123").Value = -45515.5
